$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4776.9165
$ws.Range("J17").Value = 4776.9165
$ws.Range("L17").Value = 14330.7495
$ws.Range("N17").Value = -14666.7495

$ws.Range("H33").Value = 295.42307
$ws.Range("I33").Value = 299.17648
$ws.Range("J33").Value = 288.33334
$ws.Range("K33").Value = 299.17648
$ws.Range("L33").Value = 288.33334
$ws.Range("M33").Value = -70.17648000000003
$ws.Range("N33").Value = -746.33334

$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 5000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 15000
$ws.Range("N52").Value = -15320
$ws.Range("M52").ClearContents()

$ws.Range("H70").Value = 1372510.6
$ws.Range("I70").Value = 2253.8
$ws.Range("J70").Value = 1943451
$ws.Range("K70").Value = 6761.400000000001
$ws.Range("L70").Value = 5830353
$ws.Range("M70").Value = -6491.400000000001
$ws.Range("N70").Value = -5830893

$ws.Range("H73").Value = 1372510.6
$ws.Range("I73").Value = 2253.8
$ws.Range("J73").Value = 1943451
$ws.Range("K73").Value = 6761.400000000001
$ws.Range("L73").Value = 5830353
$ws.Range("M73").Value = -5825.400000000001
$ws.Range("N73").Value = -5832225

$ws.Range("H76").Value = 5005647
$ws.Range("J76").Value = 9250.5
$ws.Range("L76").Value = 9250.5
$ws.Range("N76").Value = -9880.5

$ws.Range("H79").Value = 5005647
$ws.Range("J79").Value = 9250.5
$ws.Range("L79").Value = 9250.5
$ws.Range("N79").Value = -11434.5

$ws.Range("H100").Value = 5890.875
$ws.Range("I100").Value = 2053.6365
$ws.Range("K100").Value = 2053.6365
$ws.Range("M100").Value = -1512.6365

$ws.Range("H112").Value = 2422.95
$ws.Range("J112").Value = 2422.95
$ws.Range("L112").Value = 7268.849999999999
$ws.Range("N112").Value = -9484.849999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36742.332
$ws.Range("I32").Value = 43695.418
$ws.Range("J32").Value = 9799.125
$ws.Range("K32").Value = 43695.418
$ws.Range("L32").Value = 9799.125
$ws.Range("M32").Value = -43408.418
$ws.Range("N32").Value = -10373.125

$ws.Range("H102").Value = 1838.3334
$ws.Range("I102").Value = 1724.4375
$ws.Range("K102").Value = 1724.4375
$ws.Range("M102").Value = -102.4375

$ws.Range("H122").Value = 52604.285
$ws.Range("I122").Value = 4844.933
$ws.Range("K122").Value = 14534.799
$ws.Range("M122").Value = -12084.799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4083
$ws.Range("J20").Value = 999
$ws.Range("L20").Value = 999
$ws.Range("N20").Value = -1493

$ws.Range("H99").Value = 3527.5715
$ws.Range("I99").Value = 3846.7896
$ws.Range("J99").Value = 495
$ws.Range("K99").Value = 3846.7896
$ws.Range("L99").Value = 495
$ws.Range("M99").Value = -2348.7896
$ws.Range("N99").Value = -3491

$ws.Range("H102").Value = 10752.125
$ws.Range("I102").Value = 8002.4287
$ws.Range("J102").Value = 30000
$ws.Range("K102").Value = 8002.4287
$ws.Range("L102").Value = 30000
$ws.Range("M102").Value = -4757.4287
$ws.Range("N102").Value = -36490

$ws.Range("H105").Value = 111141480
$ws.Range("I105").Value = 166710050
$ws.Range("K105").Value = 166710050
$ws.Range("M105").Value = -166708303

$ws.Range("H107").Value = 1756.5
$ws.Range("I107").Value = 1513.5
$ws.Range("K107").Value = 1513.5
$ws.Range("M107").Value = 406.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3918.3333
$ws.Range("I58").Value = 1104.0952
$ws.Range("K58").Value = 1104.0952
$ws.Range("M58").Value = -901.0952

$ws.Range("H96").Value = 5618.5
$ws.Range("J96").Value = 5618.5
$ws.Range("L96").Value = 5618.5
$ws.Range("N96").Value = -11110.5

$ws.Range("H105").Value = 1402.8
$ws.Range("I105").Value = 1417.2858
$ws.Range("K105").Value = 1417.2858
$ws.Range("M105").Value = 329.7141999999999

$ws.Range("H122").Value = 113057.664
$ws.Range("I122").Value = 1559.5
$ws.Range("J122").Value = 144914.28
$ws.Range("K122").Value = 4678.5
$ws.Range("L122").Value = 434742.84
$ws.Range("M122").Value = -2228.5
$ws.Range("N122").Value = -439642.84

$ws.Range("H132").Value = 70587
$ws.Range("I132").Value = 8058.1665
$ws.Range("K132").Value = 24174.4995
$ws.Range("M132").Value = -21644.4995

$ws.Range("H134").Value = 1754.2273
$ws.Range("I134").Value = 1298.8524
$ws.Range("K134").Value = 3896.5572
$ws.Range("M134").Value = -1361.5572

$ws.Range("H136").Value = 3918.3333
$ws.Range("I136").Value = 1104.0952
$ws.Range("K136").Value = 3312.2856
$ws.Range("M136").Value = -762.2856000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4613
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 4919.5
$ws.Range("K64").Value = 12000
$ws.Range("L64").Value = 14758.5
$ws.Range("M64").Value = -11730
$ws.Range("N64").Value = -15298.5

$ws.Range("H67").Value = 4613
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 4919.5
$ws.Range("K67").Value = 12000
$ws.Range("L67").Value = 14758.5
$ws.Range("M67").Value = -11064
$ws.Range("N67").Value = -16630.5

$ws.Range("H68").Value = 2440
$ws.Range("I68").Value = 880
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2640
$ws.Range("L68").Value = 12000
$ws.Range("M68").Value = -1829
$ws.Range("N68").Value = -13622

$ws.Range("H71").Value = 2440
$ws.Range("I71").Value = 880
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 7920
$ws.Range("L71").Value = 36000
$ws.Range("M71").Value = -3864
$ws.Range("N71").Value = -44112

$ws.Range("H137").Value = 8936.875
$ws.Range("I137").Value = 12132.333
$ws.Range("J137").Value = 4828.4287
$ws.Range("K137").Value = 36396.999
$ws.Range("L137").Value = 14485.2861
$ws.Range("M137").Value = -31296.999
$ws.Range("N137").Value = -24685.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 31697.334
$ws.Range("J53").Value = 33786.25
$ws.Range("L53").Value = 33786.25
$ws.Range("N53").Value = -35048.25

$ws.Range("H107").Value = 1628.2307
$ws.Range("J107").Value = 1344.8334
$ws.Range("L107").Value = 1344.8334
$ws.Range("N107").Value = -5184.8334

$ws.Range("H122").Value = 2475.5715
$ws.Range("J122").Value = 2777.5
$ws.Range("L122").Value = 8332.5
$ws.Range("N122").Value = -13232.5

$ws.Range("H132").Value = 3205.5305
$ws.Range("I132").Value = 2346.5
$ws.Range("K132").Value = 7039.5
$ws.Range("M132").Value = -4509.5

$ws.Range("H134").Value = 77038.17999999999
$ws.Range("J134").Value = 77038.17999999999
$ws.Range("L134").Value = 231114.54
$ws.Range("N134").Value = -236184.54

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6569.931
$ws.Range("J46").Value = 8287.091
$ws.Range("L46").Value = 8287.091
$ws.Range("N46").Value = -8663.091

$ws.Range("H55").Value = 550.9375
$ws.Range("I55").Value = 325
$ws.Range("J55").Value = 841.4286
$ws.Range("K55").Value = 325
$ws.Range("L55").Value = 841.4286
$ws.Range("M55").Value = -152
$ws.Range("N55").Value = -1187.4286

$ws.Range("H100").Value = 22729752
$ws.Range("I100").Value = 41667880
$ws.Range("J100").Value = 3997.6
$ws.Range("K100").Value = 41667880
$ws.Range("L100").Value = 3997.6
$ws.Range("M100").Value = -41667339
$ws.Range("N100").Value = -5079.6

$ws.Range("H122").Value = 27837298
$ws.Range("I122").Value = 41670570
$ws.Range("J122").Value = 170759.17
$ws.Range("K122").Value = 125011710
$ws.Range("L122").Value = 512277.51
$ws.Range("M122").Value = -125009260
$ws.Range("N122").Value = -517177.51

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 70330
$ws.Range("J94").Value = 70330
$ws.Range("L94").Value = 70330
$ws.Range("N94").Value = -72132

$ws.Range("H112").Value = 24997.5
$ws.Range("J112").Value = 24997.5
$ws.Range("L112").Value = 24997.5
$ws.Range("N112").Value = -27951.5

$ws.Range("H123").Value = 49990
$ws.Range("J123").Value = 49990
$ws.Range("L123").Value = 49990
$ws.Range("N123").Value = -59790

$ws.Range("H132").Value = 6012.107
$ws.Range("I132").Value = 5054.7334
$ws.Range("J132").Value = 7116.769
$ws.Range("K132").Value = 15164.2002
$ws.Range("L132").Value = 21350.307
$ws.Range("M132").Value = -12634.2002
$ws.Range("N132").Value = -26410.307
